# Comment out the cgi/command rows that update the search directory's
# index.cgi on both the CAT and CIM worksheets.

$wb = $excel.ActiveWorkbook

# --- CAT sheet (rows 22/23 hold the "search" cgi + chmod command) ---
$catSheet = $wb.Worksheets.Item("CAT")
$catSheet.Range("A22").Value = "*cgi"
$catSheet.Range("A23").Value = "*command"
$catSheet.Activate()
$catSheet.Range("A22:A23").Select()

# --- CIM sheet (rows 14/15 hold the equivalent "search" cgi + command) ---
$cimSheet = $wb.Worksheets.Item("CIM")
$cimSheet.Range("A14").Value = "*cgi"
$cimSheet.Range("A15").Value = "*command"
$cimSheet.Activate()
$cimSheet.Range("B40").Select()

# Re-activate CAT sheet/selection to match the saved view state.
$catSheet.Activate()
$catSheet.Range("A22:A23").Select()
